$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C comments (results) for rows 1-17
$ws.Range("C1").Value  = "ok logando com sucesso"
$ws.Range("C2").Value  = "não realiza o loguin conforme esperado"
$ws.Range("C3").Value  = "não realiza o loguin conforme esperado"
$ws.Range("C4").Value  = "funcionando corretamente"
$ws.Range("C5").Value  = "não é autorizado conforme esperado"
$ws.Range("C6").Value  = "funcionando corretamente"
$ws.Range("C7").Value  = "não realiza a retirada, mas não informa isso ao usuario, dando a entender que realizou"
$ws.Range("C8").Value  = "funcionando corretamente"
$ws.Range("C9").Value  = "mesma situação de quando se tenta fazer saida de equipamento deslogado"
$ws.Range("C10").Value = "funcionando corretamente"
$ws.Range("C11").Value = "informa que não é autorizado"
$ws.Range("C12").Value = "funcionando corretamente"
$ws.Range("C13").Value = "ok logando com sucesso"
$ws.Range("C14").Value = "ok logando com sucesso"
$ws.Range("C15").Value = "ok funcionando corretamente"
$ws.Range("C16").Value = "ok funcionando corretamente"
$ws.Range("C17").Value = "ok delogando corretamente"

# C3's font color becomes "automatic" (matches the no-color style used elsewhere, e.g. C2)
$ws.Range("C3").Font.ColorIndex = -4105

# New row 18 data
$ws.Range("B18").Value = "logando corretamente?"
$ws.Range("C18").Value = "ok logando com sucesso"

# Update the active selection to D10
$ws.Range("D10").Select()
